$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Step 1: split paragraph 3 ("Я вот сижу и ем свою лапшу,") into two runs:
#   "Я вот сижу и ем свою " + "лапшу,"
# Both fragments must keep the explicit white color (w:val="FFFFFF" with
# w:themeColor="background1"). Paragraph 2 already contains a run
# (" лапшу,") carrying exactly that rPr, so we borrow its FormattedText
# (without ever touching the source) and drop it onto the split point in
# paragraph 3 - this yields a genuine two-run split with correct formatting.
# ---------------------------------------------------------------------------
$p2Range = $d.Paragraphs(2).Range
$p2Search = $d.Range($p2Range.Start, $p2Range.End)
$p2Search.Find.Execute("лапшу,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$sourceFormatted = $d.Range($p2Search.Start, $p2Search.End).FormattedText

$p3Range = $d.Paragraphs(3).Range
$p3Search = $d.Range($p3Range.Start, $p3Range.End)
$p3Search.Find.Execute("лапшу,", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitDest = $d.Range($p3Search.Start, $p3Search.End)
$splitDest.FormattedText = $sourceFormatted

# ---------------------------------------------------------------------------
# Step 2: insert the two new lyric lines as brand-new paragraphs right after
# paragraph 3.
# ---------------------------------------------------------------------------
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()
$newPara1 = $d.Paragraphs(4)
$newPara1.Range.Text = "А лучше бы сидел и пел гачу."

$newPara1.Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs(5)
$newPara2.Range.Text = "И я опять к нему-то и лечу,"

# ---------------------------------------------------------------------------
# Step 3: the original 4th paragraph (now the 6th) changes its wording from
# "А лучше бы сидел и пел гачу." to "Скажи, зачем, куда, и для чего?" while
# keeping its existing two-run split.
# ---------------------------------------------------------------------------
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Find.Execute("А лучше бы сидел и пел ", $true, $false, $false, $false, $false, $true, 1, $false, "Скажи, зачем, куда, и для ", 2)

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastRange = $d.Range($lastPara.Range.Start, $lastPara.Range.End)
$lastRange.Find.Execute("гачу.", $true, $false, $false, $false, $false, $true, 1, $false, "чего?", 2)
